{"js": "// Remove the redundant trailing/leading spaces that were introduced around the\n// ${optionBlock} / ${optionSelected} merge-field placeholders.\n//\n//   \"${optionBlock} \"                                  -> \"${optionBlock}\"\n//   \"${optionNumber}) ${optionText} ${optionSelected} \" -> \"${optionNumber}) ${optionText}${optionSelected}\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph containing \"${optionBlock}\" ---------------------------------\n// Drop the trailing run that holds nothing but a single space character.\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text\");\n}\nawait context.sync();\n\nconst optionBlockParagraph = paragraphs.items.find(\n  (p) => p.text === \"${optionBlock} \"\n);\nif (optionBlockParagraph) {\n  const trailingSpace = optionBlockParagraph.search(\" \", { matchCase: true });\n  trailingSpace.load(\"items\");\n  await context.sync();\n  if (trailingSpace.items.length > 0) {\n    trailingSpace.items[trailingSpace.items.length - 1].delete();\n    await context.sync();\n  }\n}\n\n// --- Paragraph containing \"${optionNumber}) ${optionText} ${optionSelected} \" ---\nconst optionLineParagraph = paragraphs.items.find((p) =>\n  p.text.indexOf(\"${optionSelected}\") !== -1\n);\nif (optionLineParagraph) {\n  // Strip the leading space glued to \"${optionSelected}\".\n  const selectedWithSpace = optionLineParagraph.search(\" ${optionSelected}\", {\n    matchCase: true,\n  });\n  selectedWithSpace.load(\"items\");\n  await context.sync();\n  if (selectedWithSpace.items.length > 0) {\n    selectedWithSpace.items[0].insertText(\n      \"${optionSelected}\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n\n  // Remove the now-orphaned trailing run that is just a single space.\n  optionLineParagraph.load(\"text\");\n  await context.sync();\n  const trailingSpace2 = optionLineParagraph.search(\" \", { matchCase: true });\n  trailingSpace2.load(\"items\");\n  await context.sync();\n  if (trailingSpace2.items.length > 0) {\n    trailingSpace2.items[trailingSpace2.items.length - 1].delete();\n    await context.sync();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the redundant trailing/leading spaces that were introduced around the\n# ${optionBlock} / ${optionSelected} merge-field placeholders.\n#\n#   \"${optionBlock} \"                                   -> \"${optionBlock}\"\n#   \"${optionNumber}) ${optionText} ${optionSelected} \"  -> \"${optionNumber}) ${optionText}${optionSelected}\"\n\n$d = $word.ActiveDocument\n\n# --- \"${optionBlock} \" -> \"${optionBlock}\" ---------------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchWildcards = $false\nif ($rng.Find.Execute(\"`${optionBlock} \")) {\n    $rng.Text = \"`${optionBlock}\"\n}\n\n# --- drop the space glued in front of \"${optionSelected}\" ------------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.MatchWildcards = $false\nif ($rng2.Find.Execute(\" `${optionSelected}\")) {\n    $rng2.Text = \"`${optionSelected}\"\n}\n\n# --- drop the now-orphaned trailing space after \"${optionSelected}\" --------\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.MatchWildcards = $false\nif ($rng3.Find.Execute(\"`${optionSelected} \")) {\n    $rng3.Text = \"`${optionSelected}\"\n}\n"}
